# Update "想去人数" (want-to-go count) values in column F across the four sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 321
$ws1.Range("F9").Value  = 85
$ws1.Range("F16").Value = 1049
$ws1.Range("F17").Value = 1015
$ws1.Range("F19").Value = 1662
$ws1.Range("F20").Value = 341
$ws1.Range("F21").Value = 6065
$ws1.Range("F23").Value = 1003
$ws1.Range("F24").Value = 1005
$ws1.Range("F25").Value = 1004
$ws1.Range("F26").Value = 4210
$ws1.Range("F27").Value = 4359
$ws1.Range("F29").Value = 113
$ws1.Range("F30").Value = 1047
$ws1.Range("F31").Value = 272
$ws1.Range("F34").Value = 1018
$ws1.Range("F35").Value = 88
$ws1.Range("F37").Value = 422
$ws1.Range("F39").Value = 193
$ws1.Range("F40").Value = 43
$ws1.Range("F42").Value = 391
$ws1.Range("F44").Value = 1109
$ws1.Range("F47").Value = 3122

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 4
$ws2.Range("F9").Value  = 611
$ws2.Range("F14").Value = 250
$ws2.Range("F27").Value = 6270
$ws2.Range("F33").Value = 29

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F8").Value  = 3031
$ws3.Range("F12").Value = 526

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 3031
$ws4.Range("F7").Value  = 85
$ws4.Range("F8").Value  = 526
$ws4.Range("F14").Value = 611
$ws4.Range("F17").Value = 1049
$ws4.Range("F18").Value = 1015
$ws4.Range("F21").Value = 341
$ws4.Range("F22").Value = 6065
$ws4.Range("F24").Value = 1003
$ws4.Range("F25").Value = 1005
$ws4.Range("F26").Value = 1004
$ws4.Range("F27").Value = 4210
$ws4.Range("F28").Value = 4359
$ws4.Range("F30").Value = 113
$ws4.Range("F31").Value = 1047
$ws4.Range("F32").Value = 272
$ws4.Range("F35").Value = 1018
$ws4.Range("F37").Value = 422
$ws4.Range("F38").Value = 193
$ws4.Range("F41").Value = 391
$ws4.Range("F45").Value = 3122
$ws4.Range("F47").Value = 6270
